# Update the vm_pu results table for the 380 kV case.
# Columns B:F, I:N (column G stays the slack bus at 1 p.u., column H is
# unused/blank in this sheet) are refreshed for every data row (2-25) with
# the newly computed per-unit voltage magnitudes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$newValues = @(
    @(1.02,1.023144274002821,1.027948343873756,1.048933383778735,1.053312928831357,1.031846365469549,1.028326142833222,1.030766586132215,1.05169184797233,1.056059243142115,1.029786483906703),
    @(1.02,1.02406641753282,1.028611407756399,1.050019366910669,1.054426690783179,1.032022223477846,1.02888678057092,1.031237873599795,1.052589268111795,1.056985254471216,1.030347917814294),
    @(1.02,1.024663107842104,1.029040107192182,1.050722910726546,1.055147912219719,1.032134177949509,1.029248952418397,1.03154180110921,1.053170204136532,1.057584399103064,1.030710603987337),
    @(1.02,1.024913956171414,1.029220248054918,1.051018881344826,1.055451243322834,1.032180803127742,1.029401065490792,1.031669325487313,1.053414488019574,1.05783626792128,1.030862933077741),
    @(1.02,1.024956074672977,1.029250489513802,1.05106858787943,1.055502181539247,1.032188605852745,1.029426597480541,1.031690722896419,1.053455507772508,1.057878557112459,1.030888501325845),
    @(1.02,1.02466645969173,1.029042514576351,1.05072686471184,1.055151964835958,1.032134802688976,1.029250985527376,1.031543506066858,1.053173468040747,1.057587764632598,1.030712639983563),
    @(1.02,1.023455915934877,1.028172500564558,1.049300222715544,1.053689217242564,1.031906177597411,1.028515736406392,1.0309260720566,1.051995084350573,1.056372202404226,1.029976346724493),
    @(1.02,1.021322827432473,1.026636820088811,1.046792745524458,1.051115843927727,1.03148926469386,1.027215584213505,1.029830247372485,1.049920513317484,1.054229887711196,1.028674348166363),
    @(1.02,1.019900837171158,1.025611354347508,1.045125445837108,1.049403087737113,1.031201913128832,1.026345799907758,1.029094485635189,1.048538754000443,1.052801472361968,1.027803328667136),
    @(1.02,1.019285123478586,1.025166933024369,1.044404521445097,1.048662118800325,1.031075261802478,1.025968466792798,1.028774666694473,1.047940745142091,1.052182907863893,1.027425459695872),
    @(1.02,1.019056423139621,1.025001798027437,1.044136892492346,1.048386990118093,1.031027883901817,1.02582820231148,1.028655687783665,1.047718663530481,1.051953137985194,1.027284996022902),
    @(1.02,1.019105479988269,1.025037222611758,1.044194292819458,1.048446001619039,1.03103806172665,1.025858294312931,1.028681217474295,1.047766298728825,1.052002424739318,1.027315130758447),
    @(1.02,1.019266218970168,1.025153284072835,1.044382396000658,1.048639374533298,1.031071352340196,1.02595687464909,1.028764835608023,1.04792238687778,1.052163915165448,1.02741385108999),
    @(1.02,1.0193652559936,1.025224785833436,1.044498313050421,1.048758531246048,1.031091819548109,1.026017599198894,1.028816331152051,1.048018563971522,1.05226341376173,1.02747466187562),
    @(1.02,1.019941700473364,1.025640841052704,1.045173312870014,1.049452277444977,1.031210271699734,1.026370827342081,1.029115685131479,1.048578448253328,1.052842523389471,1.027828391643287),
    @(1.02,1.020303293909338,1.025901718462497,1.045596997561329,1.049887624461391,1.031283977943708,1.026592208172249,1.029303133410072,1.048929729825929,1.053205769995428,1.028050086859626),
    @(1.02,1.020514206894433,1.026053846339746,1.045844224792678,1.050141619240571,1.031326754719062,1.02672126724436,1.02941235016353,1.049134655558394,1.053417640638177,1.028179329210423),
    @(1.02,1.020586123036276,1.026105711618526,1.045928539698507,1.050228235786598,1.031341304039795,1.026765261394473,1.029449570104548,1.049204534875285,1.05348988220486,1.028223385837275),
    @(1.02,1.020264498181097,1.025873732656031,1.045551529950379,1.049840909201842,1.031276092178451,1.026568463197488,1.029283034251339,1.048892037623172,1.053166797583401,1.028026308164277),
    @(1.02,1.01921888525547,1.025119108433195,1.044327000037211,1.048582428263328,1.031061558299135,1.025927848122735,1.028740217222641,1.047876421551136,1.052116360464834,1.027384783342638),
    @(1.02,1.018561485336416,1.024644316260491,1.043557983227456,1.047791750987705,1.030924740397219,1.025524454086466,1.028397862903947,1.047238126756622,1.051455865770975,1.026980816440568),
    @(1.02,1.018909983553503,1.024896043471541,1.043965568833104,1.048210848980228,1.030997453114189,1.025738358882258,1.028579451979207,1.047576473900223,1.051806010527366,1.027195025005705),
    @(1.02,1.020282028295601,1.025886378353295,1.045572074500259,1.049862017618209,1.031279656080133,1.026579192741524,1.029292116562669,1.048909069017031,1.053184407554058,1.028037052945497),
    @(1.02,1.021874272667048,1.027034131466416,1.047440222649674,1.051780625904787,1.031598707678006,1.027552239436629,1.030114466636552,1.050456613821585,1.054783765262388,1.029011481478518)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range($columns[$j] + $rowNum).Value = $rowValues[$j]
    }
}
